$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text so numeric-looking values (e.g. "1.00", "0.810") are not
# auto-converted to floating point numbers by Excel, which would also drop trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '58.903.85'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '2.491.29'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '534.61'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').Value = '136.24'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').Value = '2.511.60'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').Value = '5.32'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').Value = '0.345'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '2.958.62'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '58.783.64'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '22.84'
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').Value = '2.507.39'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').Value = '11.04'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').Value = '322.63'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').Value = '65.44'
$ws.Range('E24').Value = '  +3.44%  '
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').Value = '7.52'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').Value = '6.67'
$ws.Range('E29').Value = '  -4.21%  '
$ws.Range('D30').Value = '0.0₃0764'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '166.32'
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('E33').Value = '  +3.58%  '
$ws.Range('D34').Value = '0.997'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').Value = '18.41'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = '4.07'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('E38').Value = '  -4.15%  '
$ws.Range('D39').Value = '36.64'
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = '3.59'
$ws.Range('E41').Value = '  -2.64%  '
$ws.Range('D42').Value = '282.94'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '130.75'
$ws.Range('E45').Value = '  +5.75%  '
$ws.Range('D46').Value = '0.603'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('D47').Value = '10.89'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').Value = '0.0922'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('D49').Value = '0.0504'
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('D50').Value = '0.0219'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('D51').Value = '17.14'
$ws.Range('E51').Value = '  -4.35%  '

# Restore original (default) cell styling now that text values are safely stored,
# since the source workbook did not apply any explicit style/number format to these cells.
$ws.Range("D2:E51").ClearFormats()
